$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# --- Header row (row 1), columns X:AM ---
$hdr = New-Object "object[,]" 1,16
$hdr[0,0] = "he_beta1_opt"
$hdr[0,1] = "he_beta2_opt"
$hdr[0,2] = "he_beta3_opt"
$hdr[0,3] = "he_beta4_opt"
$hdr[0,4] = "he_beta1_err"
$hdr[0,5] = "he_beta2_err"
$hdr[0,6] = "he_beta3_err"
$hdr[0,7] = "he_beta4_err"
$hdr[0,8] = "ne_beta1_opt"
$hdr[0,9] = "ne_beta2_opt"
$hdr[0,10] = "ne_beta3_opt"
$hdr[0,11] = "ne_beta4_opt"
$hdr[0,12] = "ne_beta1_err"
$hdr[0,13] = "ne_beta2_err"
$hdr[0,14] = "ne_beta3_err"
$hdr[0,15] = "ne_beta4_err"
$ws.Range("X1:AM1").Value = $hdr

# Match header formatting (bold, centered, bordered) used by existing headers
$ws.Range("W1").Copy()
$ws.Range("X1:AM1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Data rows (rows 2:15), columns X:AM ---
$data = New-Object "object[,]" 14,16
$data[0,0] = -0.08712342836338005
$data[0,1] = 1.8526327947492
$data[0,2] = 0.03903878773608511
$data[0,3] = 1.055641995004353
$data[0,4] = 0.03049438544290839
$data[0,5] = 0.05519615582899765
$data[0,6] = 0.04177393128203785
$data[0,7] = 0.0524008869153552
$data[0,8] = 0.00665910195544306
$data[0,9] = 1.033463581479136
$data[0,10] = -0.08253613735772264
$data[0,11] = 0.1150853695187685
$data[0,12] = 0.01741669925737621
$data[0,13] = 0.02593507377490701
$data[0,14] = 0.02391484108552423
$data[0,15] = 0.02734216073984901
$data[1,0] = -0.5126359163577923
$data[1,1] = 1.871606923575062
$data[1,2] = -0.5931652978255164
$data[1,3] = 1.019959238056239
$data[1,4] = 0.03574992840703234
$data[1,5] = 0.06157534226017874
$data[1,6] = 0.04828991608573741
$data[1,7] = 0.05784469344604185
$data[1,8] = -0.3446416363958164
$data[1,9] = 1.040714088141967
$data[1,10] = -0.2766355775568772
$data[1,11] = 0.1078647245766133
$data[1,12] = 0.01684739458137365
$data[1,13] = 0.0244767250262799
$data[1,14] = 0.02272331680010336
$data[1,15] = 0.02576136461521105
$data[2,0] = -0.5852403539004787
$data[2,1] = 1.871603019760021
$data[2,2] = -0.7210711192179852
$data[2,3] = 1.055928244678011
$data[2,4] = 0.03458477848570824
$data[2,5] = 0.0586307222813099
$data[2,6] = 0.04682960648179663
$data[2,7] = 0.05539533024976183
$data[2,8] = -0.4831013453718521
$data[2,9] = 1.044212186911499
$data[2,10] = -0.2141712738799632
$data[2,11] = 0.1038006018817129
$data[2,12] = 0.01680896756964631
$data[2,13] = 0.02385881327048233
$data[2,14] = 0.02205300049460752
$data[2,15] = 0.02509042756516555
$data[3,0] = -0.2763922288423692
$data[3,1] = 1.868770633037766
$data[3,2] = -0.2777965550682772
$data[3,3] = 1.03484879991801
$data[3,4] = 0.03214935443139809
$data[3,5] = 0.05754567281419946
$data[3,6] = 0.04376570468410464
$data[3,7] = 0.05422602340177658
$data[3,8] = -0.3092617594757713
$data[3,9] = 1.036180243318022
$data[3,10] = 0.02857495851809309
$data[3,11] = 0.07702434678131423
$data[3,12] = 0.01732820140550633
$data[3,13] = 0.02527703618682386
$data[3,14] = 0.02327792598306839
$data[3,15] = 0.02661525552801863
$data[4,0] = 0.1735721045488523
$data[4,1] = 1.862423122043508
$data[4,2] = 0.3755790187415791
$data[4,3] = 1.018312054376965
$data[4,4] = 0.03079794617951738
$data[4,5] = 0.05560310569781609
$data[4,6] = 0.04267507614821369
$data[4,7] = 0.05234326299672853
$data[4,8] = 0.02482769548480739
$data[4,9] = 1.031301466077639
$data[4,10] = 0.246795773099099
$data[4,11] = 0.09924168042852502
$data[4,12] = 0.01724150514849727
$data[4,13] = 0.02565916990434255
$data[4,14] = 0.02382441613209431
$data[4,15] = 0.02705518035138128
$data[5,0] = 0.4085050493843264
$data[5,1] = 1.858930065815635
$data[5,2] = 0.7487786662396047
$data[5,3] = 1.056073089578756
$data[5,4] = 0.03285344856747372
$data[5,5] = 0.05753085994693838
$data[5,6] = 0.04630014879396727
$data[5,7] = 0.05453348258475243
$data[5,8] = 0.2699357901328576
$data[5,9] = 1.032324789309688
$data[5,10] = 0.2960107337442034
$data[5,11] = 0.1126692758946162
$data[5,12] = 0.01861432408007889
$data[5,13] = 0.02726609967958301
$data[5,14] = 0.02539097595154257
$data[5,15] = 0.02875084949266605
$data[6,0] = 0.2282000581129237
$data[6,1] = 1.853389385190773
$data[6,2] = 0.4933187541564643
$data[6,3] = 1.103257446127201
$data[6,4] = 0.03238827344259639
$data[6,5] = 0.05805591825850019
$data[6,6] = 0.04518855843098584
$data[6,7] = 0.05553616075293357
$data[6,8] = 0.2272243098448672
$data[6,9] = 1.028234644944094
$data[6,10] = 0.107184209979883
$data[6,11] = 0.1262866174274177
$data[6,12] = 0.01896362925158791
$data[6,13] = 0.02788386531910085
$data[6,14] = 0.0257546060653278
$data[6,15] = 0.02943638248310589
$data[7,0] = -0.2380628298019643
$data[7,1] = 1.85566495971684
$data[7,2] = -0.1847009399468766
$data[7,3] = 1.141836410728529
$data[7,4] = 0.03369688456674522
$data[7,5] = 0.06037531510983696
$data[7,6] = 0.04583680268353071
$data[7,7] = 0.05809610736311984
$data[7,8] = -0.09836401841795078
$data[7,9] = 1.034791571697421
$data[7,10] = -0.1645654016641566
$data[7,11] = 0.1381212413143265
$data[7,12] = 0.01714420184260041
$data[7,13] = 0.02548097750883034
$data[7,14] = 0.0235466097283301
$data[7,15] = 0.02687028008469364
$data[8,0] = -0.5781401935777103
$data[8,1] = 1.844981263941911
$data[8,2] = -0.691726151668161
$data[8,3] = 1.145367279149382
$data[8,4] = 0.03821696313080116
$data[8,5] = 0.06445416824677141
$data[8,6] = 0.05160157989756536
$data[8,7] = 0.06222747810997305
$data[8,8] = -0.419336562286888
$data[8,9] = 1.038141491040083
$data[8,10] = -0.2926238010951904
$data[8,11] = 0.1303017498290538
$data[8,12] = 0.01714867455688104
$data[8,13] = 0.02460120525356856
$data[8,14] = 0.02287585792445428
$data[8,15] = 0.02591919151411351
$data[9,0] = -0.5269693751944676
$data[9,1] = 1.847351347671347
$data[9,2] = -0.6377269520337167
$data[9,3] = 1.085257725019681
$data[9,4] = 0.03651621088090163
$data[9,5] = 0.06231965843518771
$data[9,6] = 0.04947509842904489
$data[9,7] = 0.05952959307481204
$data[9,8] = -0.462018153894671
$data[9,9] = 1.043341715913945
$data[9,10] = -0.1420586944334745
$data[9,11] = 0.1239891307390348
$data[9,12] = 0.0165508607498096
$data[9,13] = 0.02358450536469075
$data[9,14] = 0.02173677413933552
$data[9,15] = 0.02481701649024841
$data[10,0] = -0.1147804842258518
$data[10,1] = 1.843809520414935
$data[10,2] = -0.06380637873147545
$data[10,3] = 1.085127806008699
$data[10,4] = 0.03145143651882459
$data[10,5] = 0.05672910246754853
$data[10,6] = 0.04304398559188468
$data[10,7] = 0.05423713680492617
$data[10,8] = -0.2041934460210911
$data[10,9] = 1.031031523057854
$data[10,10] = 0.1286061595095815
$data[10,11] = 0.09872100005180956
$data[10,12] = 0.01713338783570813
$data[10,13] = 0.02526310925369807
$data[10,14] = 0.0233343599193656
$data[10,15] = 0.02663883452639255
$data[11,0] = 0.3035329079550512
$data[11,1] = 1.837927855084798
$data[11,2] = 0.5599688671916502
$data[11,3] = 1.113440293552569
$data[11,4] = 0.03172920809552297
$data[11,5] = 0.05615558414298819
$data[11,6] = 0.04423059709132978
$data[11,7] = 0.05402198955533404
$data[11,8] = 0.1360365966151064
$data[11,9] = 1.031717401729078
$data[11,10] = 0.3017051504957342
$data[11,11] = 0.1185241082462293
$data[11,12] = 0.01881354812344437
$data[11,13] = 0.02788852393159148
$data[11,14] = 0.0259843273120923
$data[11,15] = 0.02941458858194177
$data[12,0] = 0.4152752140684409
$data[12,1] = 1.839482529496212
$data[12,2] = 0.7529715082518194
$data[12,3] = 1.15145242619626
$data[12,4] = 0.03128674644152071
$data[12,5] = 0.05445242417392039
$data[12,6] = 0.04406992321960498
$data[12,7] = 0.05269954882236389
$data[12,8] = 0.316402902342739
$data[12,9] = 1.028196864635238
$data[12,10] = 0.252098953401878
$data[12,11] = 0.1350260113252573
$data[12,12] = 0.01979804072771231
$data[12,13] = 0.02880338353162175
$data[12,14] = 0.02676967660380434
$data[12,15] = 0.03041371217956687
$data[13,0] = 0.04969147956619423
$data[13,1] = 1.834704276414616
$data[13,2] = 0.2466713392785645
$data[13,3] = 1.190103733134485
$data[13,4] = 0.0308296226927953
$data[13,5] = 0.05561268915107377
$data[13,6] = 0.04258240714610513
$data[13,7] = 0.05424771667484846
$data[13,8] = 0.1210429013756779
$data[13,9] = 1.028816923443945
$data[13,10] = -0.008729588579138418
$data[13,11] = 0.130556850691652
$data[13,12] = 0.01889719612370874
$data[13,13] = 0.02801966007137047
$data[13,14] = 0.02584143032936607
$data[13,15] = 0.02957909204378071
$ws.Range("X2:AM15").Value = $data

Write-Output "PAD fitting result export complete"
